$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "SAP number" with "Vendor number" in cell B1
$ws.Range("B1").Value = "Vendor number"

# Make B1 the active/selected cell, matching the recorded selection in the diff
$ws.Range("B1").Select()
